# Insert 4 new weekly price rows for "Vega Central Mapocho de Santiago - Zanahoria"
# right before the current row 513, shifting the existing rows (513-573) down to
# (517-577), and fill the newly inserted rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows above the current row 513 (one at a time, so each new
# row lands above the previous insertion point and the original row 513 ends
# up at row 517).
$ws.Rows.Item(513).Insert()
$ws.Rows.Item(513).Insert()
$ws.Rows.Item(513).Insert()
$ws.Rows.Item(513).Insert()

# Columns A, B, C, E, F, G, H, N, Q, R are constant for every record in this
# sheet (same market/region/category/unit metadata) - reuse the values from
# the row right below (the old row 513, now shifted to row 517).
$constA = $ws.Cells.Item(517, 1).Value()
$constB = $ws.Cells.Item(517, 2).Value()
$constC = $ws.Cells.Item(517, 3).Value()
$constE = $ws.Cells.Item(517, 5).Value()
$constF = $ws.Cells.Item(517, 6).Value()
$constG = $ws.Cells.Item(517, 7).Value()
$constH = $ws.Cells.Item(517, 8).Value()
$constN = $ws.Cells.Item(517, 14).Value()
$constQ = $ws.Cells.Item(517, 17).Value()
$constR = $ws.Cells.Item(517, 18).Value()

$newRows = @(
    @{ Row = 513; D = 44449; I = "Primera"; J = 196; K = 5000; L = 5500; M = 5250; O = "Chillán";              P = 262 },
    @{ Row = 514; D = 44449; I = "Primera"; J = 214; K = 5500; L = 6000; M = 5750; O = "Región Metropolitana"; P = 288 },
    @{ Row = 515; D = 44449; I = "Segunda"; J = 106; K = 4000; L = 4500; M = 4250; O = "Chillán";              P = 212 },
    @{ Row = 516; D = 44449; I = "Segunda"; J = 133; K = 4500; L = 5000; M = 4748; O = "Región Metropolitana"; P = 237 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $constN
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
}

Write-Output "done"
